$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in missing "add_start" dates (column C) for a few policy rows ---
# Match the existing date formatting (copy format from the neighboring
# "add_end" cell in column D, which already carries the custom date style)
# so the new cells reuse the workbook's existing date-format style instead
# of creating a duplicate one.

$ws.Range("D14").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value = 43912

$ws.Range("D15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value = 43909

$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 43912

$excel.CutCopyMode = 0

# --- Column widths (best fit on the data columns) ---
$ws.Columns.Item(1).ColumnWidth = 40.08984375
$ws.Columns.Item(2).ColumnWidth = 21.453125
$ws.Columns.Item(3).ColumnWidth = 10.08984375
$ws.Columns.Item(4).ColumnWidth = 10.08984375

# --- Updated view / selection (scrolled down, C15 selected) ---
$ws.Activate()
$ws.Range("C15").Select() | Out-Null
